# sp_index.xlsx update — add a new species row for "NANUE, NANUE PARA, PISI O PUA"
# right after the existing NANUE entry, keeping the same sc_name / en_name
# (Kiphosus sandwicensis / Pacific chub) that the adjacent row already uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new blank row at 232, pushing everything currently at/after
# row 232 (starting with "QUIMERA NEGRA") down by one.
$ws.Rows("232:232").Insert()

# Populate the new row.
$ws.Range("A232").Value = "NANUE, NANUE PARA, PISI O PUA"
$ws.Range("B232").Value = "Kiphosus sandwicensis "
$ws.Range("C232").Value = "Pacific chub"

# Match the author's final selection/view state.
$ws.Range("B232:C232").Select() | Out-Null
